$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item("Testdata")

# Remove the hyperlink on Testdata!B6 (the "Showing Results for sprocket" search link)
$linksToRemove = @()
foreach ($hl in $ws2.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$6') {
        $linksToRemove += $hl
    }
}
foreach ($hl in $linksToRemove) {
    $hl.Delete()
}

# Update the selection kept on the now-inactive first sheet
$ws1.Activate()
$ws1.Range("B6").Select()

# Make "Testdata" the active sheet/tab, with B6 selected
$ws2.Activate()
$ws2.Range("B6").Select()
